# Updated iron ARs figs and shiny app
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix "µg" unit label -> "ug" everywhere it is used (Vitamin B12, Iodine, Selenium rows) ---
$ws.Range("C9").Value = "ug"

# --- Re-sort the "Mineral" block (rows 13-17) by descending "Billions of people" (col E),
#     reflecting updated Iron estimate that now outranks Iodine/Zinc/Magnesium ---

# Row 13 -> Iron (was Iodine)
$ws.Range("B13").Value = "Iron"
$ws.Range("C13").Value = "mg"
$ws.Range("D13").Value = "EFSA"
$ws.Range("E13").Value = 4.7672493213966503
$ws.Range("F13").Value = 0.62954011329410098

# Row 14 -> Iodine (was Zinc)
$ws.Range("B14").Value = "Iodine"
$ws.Range("C14").Value = "ug"
$ws.Range("D14").Value = "IOM"
$ws.Range("E14").Value = 3.7818351536102899
$ws.Range("F14").Value = 0.49941103780280999

# Row 15 -> Zinc (was Magnesium)
$ws.Range("B15").Value = "Zinc"
$ws.Range("C15").Value = "mg"
$ws.Range("D15").Value = "EFSA"
$ws.Range("E15").Value = 3.4799951323483098
$ws.Range("F15").Value = 0.45955149021650998

# Row 16 -> Magnesium (was Iron)
$ws.Range("B16").Value = "Magnesium"
$ws.Range("C16").Value = "mg"
$ws.Range("D16").Value = "IOM"
$ws.Range("E16").Value = 3.3588247903917798
$ws.Range("F16").Value = 0.44355031518653398

# Row 17 -> Selenium unit label fix only (values unchanged)
$ws.Range("C17").Value = "ug"

# --- Add a new (empty, but number-formatted) row 34, mirroring columns E/F styles ---
$ws.Range("E34").NumberFormat = "0.00"
$ws.Range("F34").NumberFormat = "0.0%"
$ws.Range("E34").ClearContents()
$ws.Range("F34").ClearContents()

# --- Column width adjustments to better fit the refreshed data ---
$ws.Columns.Item(1).ColumnWidth = 6.8333333
$ws.Columns.Item(4).ColumnWidth = 8.5
$ws.Columns.Item(5).ColumnWidth = 23
$ws.Columns.Item(6).ColumnWidth = 10

# --- Update selection / active cell to match the latest editing location ---
$ws.Range("D22").Select()
